$wb = $excel.ActiveWorkbook

# Update the ffcsaGuestSalesPrice (column M) values on the Pricelist sheet
$wsPricelist = $wb.Worksheets.Item("Pricelist")

$guestPriceUpdates = @{
    "M2" = 22.21
    "M3" = 17.48
    "M4" = 12.34
    "M5" = 10.28
    "M6" = 23.6
    "M7" = 27.07
    "M8" = 48.07
    "M9" = 63.16
    "M10" = 48.02
    "M12" = 37.28
    "M13" = 15.66
    "M14" = 26.74
    "M15" = 48.84
    "M16" = 17.48
    "M17" = 25.71
    "M18" = 14.4
    "M19" = 33.42
    "M20" = 55.94
    "M21" = 11.31
    "M22" = 24.06
    "M23" = 29.31
    "M24" = 24.06
    "M25" = 29.41
    "M26" = 15.42
    "M27" = 22.06
    "M28" = 45.12
    "M29" = 34.09
    "M31" = 19.02
    "M32" = 24.32
    "M33" = 26.86
    "M34" = 44.27
    "M35" = 74.04
    "M36" = 61.1
    "M37" = 48.96
    "M38" = 42.78
    "M39" = 33.42
    "M40" = 15.94
    "M41" = 14.7
    "M42" = 26.07
    "M43" = 12.34
    "M44" = 26.63
    "M45" = 12.34
    "M46" = 39.56
    "M47" = 26.74
    "M48" = 15.91
    "M49" = 14.7
    "M50" = 37.28
    "M51" = 14.65
    "M52" = 9.77
    "M53" = 9.77
    "M54" = 12.34
    "M55" = 27.76
    "M56" = 119.85
    "M57" = 19.92
    "M58" = 21.59
    "M59" = 50.9
    "M60" = 44.94
    "M61" = 38.05
    "M62" = 30.72
    "M63" = 66.63
    "M64" = 56.14
    "M65" = 43.19
    "M66" = 19.02
    "M67" = 12.34
    "M68" = 20.57
    "M69" = 10.28
    "M70" = 22.21
    "M71" = 43.19
    "M72" = 9.28
    "M73" = 27.08
    "M74" = 20.72
    "M75" = 31.39
    "M76" = 24.42
    "M77" = 43.96
    "M78" = 51.41
    "M80" = 22.11
    "M81" = 61.9
    "M82" = 14.91
    "M83" = 17.99
    "M84" = 20.05
    "M85" = 15.94
    "M86" = 20.05
    "M87" = 20.05
    "M88" = 22.11
    "M89" = 22.11
    "M90" = 16.97
    "M91" = 20.57
    "M92" = 22.11
    "M93" = 18.51
    "M94" = 18.51
    "M95" = 20.57
    "M96" = 18.51
    "M97" = 20.57
    "M98" = 9.25
    "M99" = 20.05
    "M100" = 20.57
    "M101" = 16.97
    "M102" = 11.31
    "M103" = 17.48
    "M104" = 20.05
    "M105" = 18.51
    "M106" = 12.34
    "M107" = 7.2
    "M108" = 20.57
    "M109" = 16.97
    "M110" = 16.97
    "M111" = 16.97
    "M112" = 17.48
    "M113" = 16.97
    "M114" = 20.57
    "M115" = 20.57
    "M116" = 20.57
    "M117" = 20.57
    "M118" = 15.42
    "M119" = 20.57
    "M120" = 20.57
    "M121" = 24.68
    "M122" = 20.57
    "M123" = 10.28
    "M124" = 15.42
    "M125" = 10.28
    "M126" = 16.45
    "M127" = 10.28
    "M128" = 22.62
    "M129" = 20.57
    "M130" = 17.07
    "M131" = 21.36
    "M132" = 22.5
    "M133" = 19.77
    "M134" = 22.5
    "M135" = 17.07
    "M136" = 21.36
    "M137" = 22.51
    "M138" = 17.07
    "M139" = 19.77
    "M140" = 22.5
    "M141" = 15.81
    "M142" = 9.49
    "M143" = 14.24
    "M144" = 23.3
    "M145" = 19.77
    "M146" = 21.36
    "M147" = 30.46
    "M148" = 33.61
    "M149" = 23.73
    "M150" = 22.15
    "M151" = 17.4
    "M152" = 15.81
    "M153" = 15.81
    "M154" = 15.81
    "M155" = 9.17
    "M156" = 18.98
    "M157" = 11.07
    "M158" = 18.98
}

foreach ($cell in $guestPriceUpdates.Keys) {
    $wsPricelist.Range($cell).Value = $guestPriceUpdates[$cell]
}

# Update the GUEST_MARKUP variable on the Variables sheet
$wsVariables = $wb.Worksheets.Item("Variables")
$wsVariables.Range("A5").Value = 0.9
